# Fuel Prod Imp Exp Balancing Priorities.xlsx - "updated 4.0 files and mdl"
#
# Data-level changes only (cosmetic/app-version/window-chrome/theme-label
# attributes are not exposed through the Excel object model and are left
# untouched):
#   - About!C1: refresh date 1/3/2024 -> 3/28/2024 (serial 45294 -> 45379)
#   - FPIEBP!B3:D3 (row 3 = "hard coal"): priorities 3,2,1 -> 1,3,2
#     (production priority 1, imports 3, exports 2)
#   - FPIEBP: active selection moves from F4 to E3

$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date shown in C1 -----------------
# (3/28/2024, stored as Excel serial date 45379)
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "FPIEBP" sheet: "hard coal" (row 3) balancing-priority reorder --------
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")
$wsFPIEBP.Activate()

$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# Move the selection/active cell to E3 (matches the saved cursor position)
$wsFPIEBP.Range("E3").Select()
